# Fixed minor defect in LaunchBrowser.java for chrome browser.
# Append new test-data rows to the "applicant" and "addBrandNewVisitor" sheets.

$wb = $excel.ActiveWorkbook

# --- applicant sheet: add two more URL rows (A4, A5) ---
$applicant = $wb.Worksheets.Item("applicant")
$applicant.Range("A4").Value = "https://www.bnitest2.sasyadev.com/web/open/tokenator?concept=connect&token=cQzzvAW7NGR1LDbnZCy95t%2FXLrJ5XPu83rMqH%2FkgAfXTEJtfi9OgWgbd0JIq%2BI3XHg87mKEqLKcLzbo8jZSnCK8sEvlvXa5lEPzTg8cxW2jWPGJqzdO3Mrto2KY7mQ8IvcQvkBIdUf8z1WWcRM7G%2F%2FuezrWLz7X3wlbgPOlxXLrt2Jl56ProgQ%3D%3D&chapterId=816&step=applicant"
$applicant.Range("A5").Value = "https://www.bnitest2.sasyadev.com/web/open/tokenator?concept=connect&token=xMCwr0RoBRUpHSpIzcwMIj6oOcvftDDdBJkJwDjIK1%2BKErqvrr3hTf1z9KNylv5NRWwWYvmilW%2FZU%2BwMMzFq7vfki8M7rPEURI2SYLZtg%2FdX%2Fhkzoa%2FTcxIiH6nrScP%2B5skm1sivwpMNdKZn%2FINqgtJhCCv%2FeJ1ZEbYRaK6vBJDEDf9C5R2kPQ%3D%3D&chapterId=816&step=applicant"

# --- addBrandNewVisitor sheet: add one more email/lastname row (A2, B2) ---
$addBrandNewVisitor = $wb.Worksheets.Item("addBrandNewVisitor")
$addBrandNewVisitor.Range("A2").Value = "dbselenium+v20200213121448@gmail.com"
$addBrandNewVisitor.Range("B2").Value = "selenium+v20200213121448"
